# Updated symbol list on Fri Dec 16 13:53:48 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, $text) {
    # Force the cell to keep a literal text value even when the text looks
    # like a number (leading apostrophe = "enter as text"), then strip the
    # quote-prefix cell style back off so no stray formatting is introduced.
    $ws.Range($range).Value = "'" + $text
    $ws.Range($range).Style = "Normal"
}

Set-TextValue "D2" "249.36"
Set-TextValue "D3" "23.87"
Set-TextValue "D5" "0.05920"
Set-TextValue "D6" "3.428"
Set-TextValue "D7" "6.575"
Set-TextValue "D9" "0.7958"
Set-TextValue "D10" "0.1484"
Set-TextValue "D11" "0.07883"
Set-TextValue "D12" "0.03331"
Set-TextValue "D13" "0.03032"
Set-TextValue "D14" "0.09254"
Set-TextValue "D15" "3.558"
Set-TextValue "D16" "0.001665"
Set-TextValue "D17" "0.04767"
Set-TextValue "D18" "0.0006075"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue "D19" "0.006228"
Set-TextValue "D20" "0.005676"
Set-TextValue "D23" "3.705"
Set-TextValue "D24" "2.211"
Set-TextValue "D25" "0.3332"
Set-TextValue "D26" "0.1253"
Set-TextValue "D27" "0.0006479"
Set-TextValue "D40" "0.04418"
Set-TextValue "D41" "0.007000"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003604"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1065"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Set-TextValue "D44" "0.009193"
Set-TextValue "D45" "0.002462"
Set-TextValue "D46" "0.00005892"
Set-TextValue "D48" "0.9907"
Set-TextValue "D49" "0.1110"
$ws.Range("E49").Value = "48BOLOBOLO"
Set-TextValue "D50" "0.00002102"
